# Update "想去人数" (column F) counts across the four sheets to reflect
# the regenerated output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 55
$ws1.Range("F3").Value  = 187
$ws1.Range("F5").Value  = 363
$ws1.Range("F7").Value  = 162
$ws1.Range("F12").Value = 172
$ws1.Range("F14").Value = 5989
$ws1.Range("F16").Value = 2311
$ws1.Range("F19").Value = 460
$ws1.Range("F20").Value = 9019
$ws1.Range("F22").Value = 2433
$ws1.Range("F25").Value = 2413
$ws1.Range("F27").Value = 233
$ws1.Range("F28").Value = 1942
$ws1.Range("F42").Value = 1512
$ws1.Range("F43").Value = 2465
$ws1.Range("F45").Value = 910
$ws1.Range("F47").Value = 1247
$ws1.Range("F48").Value = 12

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F22").Value = 42
$ws2.Range("F23").Value = 42

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 884
$ws3.Range("F4").Value = 99

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 55
$ws4.Range("F4").Value  = 884
$ws4.Range("F5").Value  = 99
$ws4.Range("F6").Value  = 363
$ws4.Range("F11").Value = 162
$ws4.Range("F16").Value = 172
$ws4.Range("F19").Value = 5989
$ws4.Range("F21").Value = 2311
$ws4.Range("F23").Value = 460
$ws4.Range("F24").Value = 9019
$ws4.Range("F27").Value = 2433
$ws4.Range("F29").Value = 2413
$ws4.Range("F31").Value = 233
$ws4.Range("F32").Value = 1942
$ws4.Range("F43").Value = 1512
$ws4.Range("F44").Value = 2465
$ws4.Range("F45").Value = 910
$ws4.Range("F50").Value = 1247
$ws4.Range("F51").Value = 42
